$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-9) for columns D, M, N, O, P, Q, S have been shuffled
# between rows according to a fixed permutation (row <- source row):
#   2<-8, 3<-5, 4<-6, 5<-9, 6<-7, 7<-2, 8<-3, 9<-4
# Capture the original values first, then write them back in their new
# positions so the shuffle is applied consistently (no overwritten sources).

$cols = @("D", "M", "N", "O", "P", "Q", "S")
$mapping = @{ 2 = 8; 3 = 5; 4 = 6; 5 = 9; 6 = 7; 7 = 2; 8 = 3; 9 = 4 }

$original = @{}
foreach ($row in 2..9) {
    $original[$row] = @{}
    foreach ($col in $cols) {
        $original[$row][$col] = $ws.Range("$col$row").Value2
    }
}

foreach ($row in 2..9) {
    $srcRow = $mapping[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $original[$srcRow][$col]
    }
}
